$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 4099.6665
$ws.Range("I58").Value = 1049.5
$ws.Range("J58").Value = 5624.75
$ws.Range("K58").Value = 3148.5
$ws.Range("L58").Value = 16874.25
$ws.Range("M58").Value = -2998.5
$ws.Range("N58").Value = -17174.25

$ws.Range("H61").Value = 555
$ws.Range("I61").Value = 555
$ws.Range("K61").Value = 1665
$ws.Range("M61").Value = -1493

$ws.Range("H80").Value = 3633070
$ws.Range("I80").Value = 8696111
$ws.Range("J80").Value = 16612.285
$ws.Range("K80").Value = 26088333
$ws.Range("L80").Value = 49836.855
$ws.Range("M80").Value = -26087335
$ws.Range("N80").Value = -51832.855

$ws.Range("H82").Value = 6719.6665
$ws.Range("I82").Value = 160
$ws.Range("K82").Value = 480
$ws.Range("M82").Value = -74

$ws.Range("H83").Value = 3633070
$ws.Range("I83").Value = 8696111
$ws.Range("J83").Value = 16612.285
$ws.Range("K83").Value = 78264999
$ws.Range("L83").Value = 149510.565
$ws.Range("M83").Value = -78260007
$ws.Range("N83").Value = -159494.565

$ws.Range("H85").Value = 6719.6665
$ws.Range("I85").Value = 160
$ws.Range("K85").Value = 480
$ws.Range("M85").Value = 924

$ws.Range("H86").Value = 2337.2307
$ws.Range("I86").Value = 2040.6666
$ws.Range("J86").Value = 3004.5
$ws.Range("K86").Value = 2040.6666
$ws.Range("L86").Value = 3004.5
$ws.Range("M86").Value = -917.6666
$ws.Range("N86").Value = -5250.5

$ws.Range("H89").Value = 2337.2307
$ws.Range("I89").Value = 2040.6666
$ws.Range("J89").Value = 3004.5
$ws.Range("K89").Value = 10203.333
$ws.Range("L89").Value = 15022.5
$ws.Range("M89").Value = -4587.333000000001
$ws.Range("N89").Value = -26254.5

$ws.Range("H99").Value = 268.83334
$ws.Range("I99").Value = 268.83334
$ws.Range("K99").Value = 806.5000200000001
$ws.Range("M99").Value = 691.4999799999999

$ws.Range("H101").Value = 233
$ws.Range("I101").Value = 244
$ws.Range("K101").Value = 732
$ws.Range("M101").Value = 890

$ws.Range("H115").Value = 724.7273
$ws.Range("I115").Value = 724.7273
$ws.Range("K115").Value = 2174.1819
$ws.Range("M115").Value = -607.1819

$ws.Range("H116").Value = 13668.4
$ws.Range("I116").Value = 15006.577
$ws.Range("J116").Value = 9802.556
$ws.Range("K116").Value = 15006.577
$ws.Range("L116").Value = 9802.556
$ws.Range("M116").Value = -11564.577
$ws.Range("N116").Value = -16686.556

$ws.Range("H118").Value = 572.1429000000001
$ws.Range("I118").Value = 572.1429000000001
$ws.Range("K118").Value = 1716.4287
$ws.Range("M118").Value = -59.42870000000016

$ws.Range("H137").Value = 30053.545
$ws.Range("J137").Value = 3319.2
$ws.Range("L137").Value = 9957.599999999999
$ws.Range("N137").Value = -15057.6

$ws.Range("H138").Value = 44146.293
$ws.Range("I138").Value = 2349.9443
$ws.Range("K138").Value = 7049.8329
$ws.Range("M138").Value = -1909.8329

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 236074
$ws.Range("I74").Value = 286623.66
$ws.Range("K74").Value = 286623.66
$ws.Range("M74").Value = -285749.66

$ws.Range("H77").Value = 236074
$ws.Range("I77").Value = 286623.66
$ws.Range("K77").Value = 1433118.3
$ws.Range("M77").Value = -1428750.3

$ws.Range("H97").Value = 882.5806
$ws.Range("I97").Value = 693.1852
$ws.Range("K97").Value = 693.1852
$ws.Range("M97").Value = -197.1852

$ws.Range("H122").Value = 3114.8708
$ws.Range("I122").Value = 3041
$ws.Range("K122").Value = 9123
$ws.Range("M122").Value = -6673

$ws.Range("H132").Value = 1646.2106
$ws.Range("I132").Value = 1292.375
$ws.Range("K132").Value = 3877.125
$ws.Range("M132").Value = -1347.125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H13").Value = 20000
$ws.Range("J13").Value = 20000
$ws.Range("L13").Value = 20000
$ws.Range("N13").Value = -20336

$ws.Range("H81").Value = 49997
$ws.Range("J81").Value = 49997
$ws.Range("L81").Value = 49997
$ws.Range("N81").Value = -52119

$ws.Range("H84").Value = 49997
$ws.Range("J84").Value = 49997
$ws.Range("L84").Value = 149991
$ws.Range("N84").Value = -160599

$ws.Range("H96").Value = 21162.5
$ws.Range("I96").Value = 9883.333000000001
$ws.Range("K96").Value = 9883.333000000001
$ws.Range("M96").Value = -7137.333000000001

$ws.Range("H134").Value = 9166.194
$ws.Range("I134").Value = 13867.6
$ws.Range("J134").Value = 3289.4375
$ws.Range("K134").Value = 41602.8
$ws.Range("L134").Value = 9868.3125
$ws.Range("M134").Value = -39067.8
$ws.Range("N134").Value = -14938.3125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1523.1666
$ws.Range("I58").Value = 1434.3636
$ws.Range("J58").Value = 2500
$ws.Range("K58").Value = 1434.3636
$ws.Range("L58").Value = 2500
$ws.Range("M58").Value = -1231.3636
$ws.Range("N58").Value = -2906

$ws.Range("H122").Value = 2372.3635
$ws.Range("I122").Value = 2210.25
$ws.Range("K122").Value = 6630.75
$ws.Range("M122").Value = -4180.75

$ws.Range("H132").Value = 78537.62
$ws.Range("I132").Value = 101049
$ws.Range("J132").Value = 3499.6667
$ws.Range("K132").Value = 303147
$ws.Range("L132").Value = 10499.0001
$ws.Range("M132").Value = -300617
$ws.Range("N132").Value = -15559.0001

$ws.Range("H134").Value = 2084.2
$ws.Range("I134").Value = 2017.375
$ws.Range("J134").Value = 3688
$ws.Range("K134").Value = 6052.125
$ws.Range("L134").Value = 11064
$ws.Range("M134").Value = -3517.125
$ws.Range("N134").Value = -16134

$ws.Range("H136").Value = 1523.1666
$ws.Range("I136").Value = 1434.3636
$ws.Range("J136").Value = 2500
$ws.Range("K136").Value = 4303.0908
$ws.Range("L136").Value = 7500
$ws.Range("M136").Value = -1753.0908
$ws.Range("N136").Value = -12600

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H127").Value = 3257.25
$ws.Range("J127").Value = 3999.6667
$ws.Range("L127").Value = 11999.0001
$ws.Range("N127").Value = -21919.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5274.25
$ws.Range("I70").Value = 4985.7
$ws.Range("J70").Value = 5562.8
$ws.Range("K70").Value = 4985.7
$ws.Range("L70").Value = 5562.8
$ws.Range("M70").Value = -4715.7
$ws.Range("N70").Value = -6102.8

$ws.Range("H73").Value = 5274.25
$ws.Range("I73").Value = 4985.7
$ws.Range("J73").Value = 5562.8
$ws.Range("K73").Value = 4985.7
$ws.Range("L73").Value = 5562.8
$ws.Range("M73").Value = -4049.7
$ws.Range("N73").Value = -7434.8

$ws.Range("H97").Value = 1052.1428
$ws.Range("I97").Value = 1052.1428
$ws.Range("K97").Value = 1052.1428
$ws.Range("M97").Value = -556.1428000000001

$ws.Range("H122").Value = 25003092
$ws.Range("I122").Value = 2594.7144
$ws.Range("K122").Value = 7784.1432
$ws.Range("M122").Value = -5334.1432

$ws.Range("H132").Value = 2462.8333
$ws.Range("I132").Value = 2327.9048
$ws.Range("K132").Value = 6983.714399999999
$ws.Range("M132").Value = -4453.714399999999

$ws.Range("H136").Value = 47999.5
$ws.Range("J136").Value = 47999.5
$ws.Range("L136").Value = 143998.5
$ws.Range("N136").Value = -149098.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H25").Value = 54949.5
$ws.Range("J25").Value = 54949.5
$ws.Range("L25").Value = 54949.5
$ws.Range("N25").Value = -55409.5

$ws.Range("H80").Value = 50000
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 50000
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 50000
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -52246

$ws.Range("H83").Value = 50000
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 50000
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 150000
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -161232

$ws.Range("H96").Value = 68996.664
$ws.Range("J96").Value = 68996.664
$ws.Range("L96").Value = 68996.664
$ws.Range("N96").Value = -74488.664

$ws.Range("H136").Value = 3015.5881
$ws.Range("I136").Value = 2442
$ws.Range("K136").Value = 7326
$ws.Range("M136").Value = -4776

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1039.0625
$ws.Range("I107").Value = 947.4545000000001
$ws.Range("J107").Value = 1240.6
$ws.Range("K107").Value = 2842.3635
$ws.Range("L107").Value = 3721.8
$ws.Range("M107").Value = -922.3635000000004
$ws.Range("N107").Value = -7561.799999999999

$ws.Range("H122").Value = 60057
$ws.Range("I122").Value = 67678.67999999999
$ws.Range("K122").Value = 203036.04
$ws.Range("M122").Value = -200586.04

$ws.Range("H136").Value = 24727.74
$ws.Range("I136").Value = 33160.95
$ws.Range("K136").Value = 99482.84999999999
$ws.Range("M136").Value = -96932.84999999999
